# Auto-generated edit script: update Shiva Profits leve-crafting values
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 561.4666999999999
$ws.Range("I11").Value = 561.4666999999999
$ws.Range("K11").Value = 561.4666999999999
$ws.Range("M11").Value = -421.4666999999999
$ws.Range("H40").Value = 2866.4167
$ws.Range("J40").Value = 4951
$ws.Range("L40").Value = 4951
$ws.Range("N40").Value = -5301
$ws.Range("H74").Value = 3064.6365
$ws.Range("I74").Value = 3064.6365
$ws.Range("K74").Value = 3064.6365
$ws.Range("M74").Value = -2128.6365
$ws.Range("H76").Value = 4900.1
$ws.Range("I76").Value = 4071.7144
$ws.Range("J76").Value = 6833
$ws.Range("K76").Value = 4071.7144
$ws.Range("L76").Value = 6833
$ws.Range("M76").Value = -3756.7144
$ws.Range("N76").Value = -7463
$ws.Range("H77").Value = 3064.6365
$ws.Range("I77").Value = 3064.6365
$ws.Range("K77").Value = 15323.1825
$ws.Range("M77").Value = -10643.1825
$ws.Range("H79").Value = 4900.1
$ws.Range("I79").Value = 4071.7144
$ws.Range("J79").Value = 6833
$ws.Range("K79").Value = 4071.7144
$ws.Range("L79").Value = 6833
$ws.Range("M79").Value = -2979.7144
$ws.Range("N79").Value = -9017
$ws.Range("H98").Value = 2032.75
$ws.Range("I98").Value = 1609
$ws.Range("J98").Value = 4999
$ws.Range("K98").Value = 1609
$ws.Range("L98").Value = 4999
$ws.Range("M98").Value = -111
$ws.Range("N98").Value = -7995
$ws.Range("H122").Value = 2032.75
$ws.Range("I122").Value = 1609
$ws.Range("J122").Value = 4999
$ws.Range("K122").Value = 4827
$ws.Range("L122").Value = 14997
$ws.Range("M122").Value = -2377
$ws.Range("N122").Value = -19897
$ws.Range("H125").Value = 1539.0526
$ws.Range("J125").Value = 2864.6667
$ws.Range("L125").Value = 25782.0003
$ws.Range("N125").Value = -30702.0003
$ws.Range("H132").Value = 6288.9854
$ws.Range("I132").Value = 3920.8447
$ws.Range("J132").Value = 20024.2
$ws.Range("K132").Value = 11762.5341
$ws.Range("L132").Value = 60072.60000000001
$ws.Range("M132").Value = -9232.534100000001
$ws.Range("N132").Value = -65132.60000000001
$ws.Range("H137").Value = 6005.0654
$ws.Range("I137").Value = 6053.3335
$ws.Range("J137").Value = 5882.5386
$ws.Range("K137").Value = 18160.0005
$ws.Range("L137").Value = 17647.6158
$ws.Range("M137").Value = -15610.0005
$ws.Range("N137").Value = -22747.6158
$ws.Range("H138").Value = 25643086
$ws.Range("I138").Value = 35715428
$ws.Range("J138").Value = 4393.636
$ws.Range("K138").Value = 107146284
$ws.Range("L138").Value = 13180.908
$ws.Range("M138").Value = -107141144
$ws.Range("N138").Value = -23460.908

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 26316724
$ws.Range("I2").Value = 35715204
$ws.Range("K2").Value = 35715204
$ws.Range("M2").Value = -35715091
$ws.Range("H32").Value = 4716.488
$ws.Range("I32").Value = 4766.074
$ws.Range("K32").Value = 4766.074
$ws.Range("M32").Value = -4479.074
$ws.Range("H116").Value = 26316724
$ws.Range("I116").Value = 35715204
$ws.Range("K116").Value = 35715204
$ws.Range("M116").Value = -35712910

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 26316724
$ws.Range("I3").Value = 35715204
$ws.Range("K3").Value = 35715204
$ws.Range("M3").Value = -35715090
$ws.Range("H94").Value = 500.5238
$ws.Range("I94").Value = 413.85715
$ws.Range("J94").Value = 673.8570999999999
$ws.Range("K94").Value = 413.85715
$ws.Range("L94").Value = 673.8570999999999
$ws.Range("M94").Value = 37.14285000000001
$ws.Range("N94").Value = -1575.8571
$ws.Range("H134").Value = 2054.215
$ws.Range("I134").Value = 1997.2106
$ws.Range("J134").Value = 3498.3333
$ws.Range("K134").Value = 5991.6318
$ws.Range("L134").Value = 10494.9999
$ws.Range("M134").Value = -3456.6318
$ws.Range("N134").Value = -15564.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2126.4062
$ws.Range("I31").Value = 1827.2174
$ws.Range("J31").Value = 2891
$ws.Range("K31").Value = 1827.2174
$ws.Range("L31").Value = 2891
$ws.Range("M31").Value = -1532.2174
$ws.Range("N31").Value = -3481
$ws.Range("H34").Value = 2126.4062
$ws.Range("I34").Value = 1827.2174
$ws.Range("J34").Value = 2891
$ws.Range("K34").Value = 1827.2174
$ws.Range("L34").Value = 2891
$ws.Range("M34").Value = -1625.2174
$ws.Range("N34").Value = -3295
$ws.Range("H58").Value = 1321.2458
$ws.Range("I58").Value = 1190.1482
$ws.Range("K58").Value = 1190.1482
$ws.Range("M58").Value = -987.1482000000001
$ws.Range("H134").Value = 2175.9194
$ws.Range("I134").Value = 2123.2075
$ws.Range("J134").Value = 2486.3333
$ws.Range("K134").Value = 6369.622499999999
$ws.Range("L134").Value = 7458.999899999999
$ws.Range("M134").Value = -3834.622499999999
$ws.Range("N134").Value = -12528.9999
$ws.Range("H136").Value = 1321.2458
$ws.Range("I136").Value = 1190.1482
$ws.Range("K136").Value = 3570.4446
$ws.Range("M136").Value = -1020.4446

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 2149.2856
$ws.Range("I22").Value = 985
$ws.Range("K22").Value = 2955
$ws.Range("M22").Value = -2786
$ws.Range("H27").Value = 2149.2856
$ws.Range("I27").Value = 985
$ws.Range("K27").Value = 2955
$ws.Range("M27").Value = -2853
$ws.Range("H38").Value = 386
$ws.Range("I38").Value = 85.07692
$ws.Range("J38").Value = 875
$ws.Range("K38").Value = 255.23076
$ws.Range("L38").Value = 2625
$ws.Range("M38").Value = 91.76924
$ws.Range("N38").Value = -3319
$ws.Range("H113").Value = 962.64514
$ws.Range("I113").Value = 795.0526
$ws.Range("J113").Value = 1228
$ws.Range("K113").Value = 2385.1578
$ws.Range("L113").Value = 3684
$ws.Range("M113").Value = -215.1578
$ws.Range("N113").Value = -8024
$ws.Range("H129").Value = 19341.076
$ws.Range("I129").Value = 762.75
$ws.Range("J129").Value = 49066.4
$ws.Range("K129").Value = 2288.25
$ws.Range("L129").Value = 147199.2
$ws.Range("M129").Value = 2711.75
$ws.Range("N129").Value = -157199.2
$ws.Range("H131").Value = 2702728.5
$ws.Range("I131").Value = 5943819
$ws.Range("J131").Value = 1820
$ws.Range("K131").Value = 17831457
$ws.Range("L131").Value = 5460
$ws.Range("M131").Value = -17826417
$ws.Range("N131").Value = -15540
$ws.Range("H132").Value = 4514.3335
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 4514.3335
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 40629.0015
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -45689.0015
$ws.Range("H139").Value = 2110.6667
$ws.Range("I139").Value = 1166.1818
$ws.Range("K139").Value = 3498.5454
$ws.Range("M139").Value = 1641.4546
$ws.Range("H140").Value = 3214.9092
$ws.Range("I140").Value = 2262.6667
$ws.Range("K140").Value = 6788.000100000001
$ws.Range("M140").Value = -1608.000100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2933.0435
$ws.Range("I132").Value = 2112.2258
$ws.Range("J132").Value = 4629.4
$ws.Range("K132").Value = 6336.6774
$ws.Range("L132").Value = 13888.2
$ws.Range("M132").Value = -3806.6774
$ws.Range("N132").Value = -18948.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 13334287
$ws.Range("I93").Value = 14286694
$ws.Range("K93").Value = 14286694
$ws.Range("M93").Value = -14285446
$ws.Range("H132").Value = 56195.7
$ws.Range("I132").Value = 58043.863
$ws.Range("J132").Value = 2599
$ws.Range("K132").Value = 174131.589
$ws.Range("L132").Value = 7797
$ws.Range("M132").Value = -171601.589
$ws.Range("N132").Value = -12857
$ws.Range("H136").Value = 11349.5
$ws.Range("I136").Value = 14236.625
$ws.Range("K136").Value = 42709.875
$ws.Range("M136").Value = -40159.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3137834.5
$ws.Range("J81").Value = 3397.4167
$ws.Range("L81").Value = 6794.8334
$ws.Range("N81").Value = -8916.8334
$ws.Range("H84").Value = 3137834.5
$ws.Range("J84").Value = 3397.4167
$ws.Range("L84").Value = 33974.167
$ws.Range("N84").Value = -44582.167
$ws.Range("H132").Value = 3437.1592
$ws.Range("I132").Value = 2754.2856
$ws.Range("J132").Value = 6092.778
$ws.Range("K132").Value = 8262.856800000001
$ws.Range("L132").Value = 18278.334
$ws.Range("M132").Value = -5732.856800000001
$ws.Range("N132").Value = -23338.334
